# Fruta / hortaliza, semanal
# Inserts a new weekly price-report group (date 2022-07-11, serial 44753)
# for "Terminal La Palmera de La Serena - Piña / Caramelo" above the
# existing data block (which starts at row 987), shifting the rest of the
# table down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before row 987 (existing rows 987:1034 shift to 991:1038)
$ws.Rows("987:990").Insert()

# Common values shared by the four new rows
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$fecha     = 44753
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100108
$producto   = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria   = "Piña"
$variedad    = "Caramelo"
$origen      = "Ecuador"

$newRows = @(
    @{ Row = 987; Calidad = "Especial"; Volumen = 216; PMin = 20000; PMax = 21000; PProm = 20500; Unidad = "`$/caja 10 unidades"; Kg = 2050; Unid = 10 },
    @{ Row = 988; Calidad = "Primera";  Volumen = 216; PMin = 20000; PMax = 21000; PProm = 20500; Unidad = "`$/caja 12 unidades"; Kg = 1708; Unid = 12 },
    @{ Row = 989; Calidad = "Segunda";  Volumen = 208; PMin = 20000; PMax = 21000; PProm = 20519; Unidad = "`$/caja 14 unidades"; Kg = 1466; Unid = 14 },
    @{ Row = 990; Calidad = "Tercera";  Volumen = 216; PMin = 20000; PMax = 21000; PProm = 20500; Unidad = "`$/caja 16 unidades"; Kg = 1281; Unid = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.PMin
    $ws.Cells.Item($row, 15).Value2 = $r.PMax
    $ws.Cells.Item($row, 16).Value2 = $r.PProm
    $ws.Cells.Item($row, 17).Value2 = $r.Unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.Kg
    $ws.Cells.Item($row, 20).Value2 = $r.Unid
}
